# List of parts updated
# STATUS column (C) for the batch of parts ordered on 2023-11-xx (rows 41-51)
# has moved from "Ordered" to "Ready".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C41").Value = "Ready"
$ws.Range("C43:C51").Value = "Ready"
